$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.023873
$ws.Cells.Item(2, 8).Value = 0.071619
$ws.Cells.Item(2, 9).Value = 0.02747901635872243
$ws.Cells.Item(2, 10).Value = 0.02747901635872243
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.569028
$ws.Cells.Item(2, 14).Value = 1.707084
$ws.Cells.Item(2, 15).Value = 0.1016535000995941
$ws.Cells.Item(2, 16).Value = 0.1016535000995941
$ws.Cells.Item(2, 17).Value = 0.013584405444
$ws.Cells.Item(2, 18).Value = 0.122259648996
$ws.Cells.Item(2, 19).Value = 0.002793338192158138
$ws.Cells.Item(2, 20).Value = 0.002793338192158137
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.023873
$ws.Cells.Item(3, 8).Value = 0.071619
$ws.Cells.Item(3, 9).Value = 0.02747901635872243
$ws.Cells.Item(3, 10).Value = 0.02747901635872243
$ws.Cells.Item(3, 15).Value = 0.1962512724671019
$ws.Cells.Item(3, 16).Value = 0.1962512724671019
$ws.Cells.Item(3, 17).Value = 0.02622592288
$ws.Cells.Item(3, 18).Value = 0.23603330592
$ws.Cells.Item(3, 19).Value = 0.005392791926543586
$ws.Cells.Item(3, 20).Value = 0.005392791926543586
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.023873
$ws.Cells.Item(4, 8).Value = 0.071619
$ws.Cells.Item(4, 9).Value = 0.02747901635872243
$ws.Cells.Item(4, 10).Value = 0.02747901635872243
$ws.Cells.Item(4, 13).Value = 3.083549
$ws.Cells.Item(4, 14).Value = 9.250647000000001
$ws.Cells.Item(4, 15).Value = 0.5508578638987945
$ws.Cells.Item(4, 16).Value = 0.5508578638987945
$ws.Cells.Item(4, 17).Value = 0.07361356527700001
$ws.Cells.Item(4, 18).Value = 0.6625220874930001
$ws.Cells.Item(4, 19).Value = 0.01513703225340587
$ws.Cells.Item(4, 20).Value = 0.01513703225340587
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.023873
$ws.Cells.Item(5, 8).Value = 0.071619
$ws.Cells.Item(5, 9).Value = 0.02747901635872243
$ws.Cells.Item(5, 10).Value = 0.02747901635872243
$ws.Cells.Item(5, 13).Value = 0.3400753333333333
$ws.Cells.Item(5, 14).Value = 1.020226
$ws.Cells.Item(5, 15).Value = 0.06075245494223394
$ws.Cells.Item(5, 16).Value = 0.06075245494223393
$ws.Cells.Item(5, 17).Value = 0.008118618432666667
$ws.Cells.Item(5, 18).Value = 0.07306756589400001
$ws.Cells.Item(5, 19).Value = 0.001669417703190193
$ws.Cells.Item(5, 20).Value = 0.001669417703190193
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.023873
$ws.Cells.Item(6, 8).Value = 0.071619
$ws.Cells.Item(6, 9).Value = 0.02747901635872243
$ws.Cells.Item(6, 10).Value = 0.02747901635872243
$ws.Cells.Item(6, 13).Value = 0.5065093333333334
$ws.Cells.Item(6, 14).Value = 1.519528
$ws.Cells.Item(6, 15).Value = 0.0904849085922755
$ws.Cells.Item(6, 16).Value = 0.09048490859227548
$ws.Cells.Item(6, 17).Value = 0.01209189731466667
$ws.Cells.Item(6, 18).Value = 0.108827075832
$ws.Cells.Item(6, 19).Value = 0.002486436283424642
$ws.Cells.Item(6, 20).Value = 0.002486436283424641
$ws.Cells.Item(7, 9).Value = 0.04107483513127341
$ws.Cells.Item(7, 10).Value = 0.04107483513127341
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.569028
$ws.Cells.Item(7, 14).Value = 1.707084
$ws.Cells.Item(7, 15).Value = 0.1016535000995941
$ws.Cells.Item(7, 16).Value = 0.1016535000995941
$ws.Cells.Item(7, 17).Value = 0.020305574504
$ws.Cells.Item(7, 18).Value = 0.182750170536
$ws.Cells.Item(7, 19).Value = 0.004175400757107712
$ws.Cells.Item(7, 20).Value = 0.004175400757107713
$ws.Cells.Item(8, 9).Value = 0.04107483513127341
$ws.Cells.Item(8, 10).Value = 0.04107483513127341
$ws.Cells.Item(8, 15).Value = 0.1962512724671019
$ws.Cells.Item(8, 16).Value = 0.1962512724671019
$ws.Cells.Item(8, 18).Value = 0.35281572672
$ws.Cells.Item(8, 19).Value = 0.008060988660888828
$ws.Cells.Item(8, 20).Value = 0.008060988660888828
$ws.Cells.Item(9, 9).Value = 0.04107483513127341
$ws.Cells.Item(9, 10).Value = 0.04107483513127341
$ws.Cells.Item(9, 13).Value = 3.083549
$ws.Cells.Item(9, 14).Value = 9.250647000000001
$ws.Cells.Item(9, 15).Value = 0.5508578638987945
$ws.Cells.Item(9, 16).Value = 0.5508578638987945
$ws.Cells.Item(9, 17).Value = 0.1100354182153333
$ws.Cells.Item(9, 18).Value = 0.990318763938
$ws.Cells.Item(9, 19).Value = 0.02262639594040843
$ws.Cells.Item(9, 20).Value = 0.02262639594040844
$ws.Cells.Item(10, 9).Value = 0.04107483513127341
$ws.Cells.Item(10, 10).Value = 0.04107483513127341
$ws.Cells.Item(10, 13).Value = 0.3400753333333333
$ws.Cells.Item(10, 14).Value = 1.020226
$ws.Cells.Item(10, 15).Value = 0.06075245494223394
$ws.Cells.Item(10, 16).Value = 0.06075245494223393
$ws.Cells.Item(10, 17).Value = 0.01213547491155556
$ws.Cells.Item(10, 18).Value = 0.109219274204
$ws.Cells.Item(10, 19).Value = 0.002495397070572376
$ws.Cells.Item(10, 20).Value = 0.002495397070572376
$ws.Cells.Item(11, 9).Value = 0.04107483513127341
$ws.Cells.Item(11, 10).Value = 0.04107483513127341
$ws.Cells.Item(11, 13).Value = 0.5065093333333334
$ws.Cells.Item(11, 14).Value = 1.519528
$ws.Cells.Item(11, 15).Value = 0.0904849085922755
$ws.Cells.Item(11, 16).Value = 0.09048490859227548
$ws.Cells.Item(11, 17).Value = 0.01807461672355555
$ws.Cells.Item(11, 18).Value = 0.162671550512
$ws.Cells.Item(11, 19).Value = 0.00371665270229606
$ws.Cells.Item(11, 20).Value = 0.00371665270229606
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.011261
$ws.Cells.Item(12, 8).Value = 0.033783
$ws.Cells.Item(12, 9).Value = 0.01296197391260307
$ws.Cells.Item(12, 10).Value = 0.01296197391260308
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.569028
$ws.Cells.Item(12, 14).Value = 1.707084
$ws.Cells.Item(12, 15).Value = 0.1016535000995941
$ws.Cells.Item(12, 16).Value = 0.1016535000995941
$ws.Cells.Item(12, 17).Value = 0.006407824308
$ws.Cells.Item(12, 18).Value = 0.057670418772
$ws.Cells.Item(12, 19).Value = 0.001317630016415733
$ws.Cells.Item(12, 20).Value = 0.001317630016415733
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.011261
$ws.Cells.Item(13, 8).Value = 0.033783
$ws.Cells.Item(13, 9).Value = 0.01296197391260307
$ws.Cells.Item(13, 10).Value = 0.01296197391260308
$ws.Cells.Item(13, 15).Value = 0.1962512724671019
$ws.Cells.Item(13, 16).Value = 0.1962512724671019
$ws.Cells.Item(13, 17).Value = 0.01237088416
$ws.Cells.Item(13, 18).Value = 0.11133795744
$ws.Cells.Item(13, 19).Value = 0.002543803874033733
$ws.Cells.Item(13, 20).Value = 0.002543803874033733
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.011261
$ws.Cells.Item(14, 8).Value = 0.033783
$ws.Cells.Item(14, 9).Value = 0.01296197391260307
$ws.Cells.Item(14, 10).Value = 0.01296197391260308
$ws.Cells.Item(14, 13).Value = 3.083549
$ws.Cells.Item(14, 14).Value = 9.250647000000001
$ws.Cells.Item(14, 15).Value = 0.5508578638987945
$ws.Cells.Item(14, 16).Value = 0.5508578638987945
$ws.Cells.Item(14, 17).Value = 0.034723845289
$ws.Cells.Item(14, 18).Value = 0.312514607601
$ws.Cells.Item(14, 19).Value = 0.00714020526140843
$ws.Cells.Item(14, 20).Value = 0.007140205261408431
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.011261
$ws.Cells.Item(15, 8).Value = 0.033783
$ws.Cells.Item(15, 9).Value = 0.01296197391260307
$ws.Cells.Item(15, 10).Value = 0.01296197391260308
$ws.Cells.Item(15, 13).Value = 0.3400753333333333
$ws.Cells.Item(15, 14).Value = 1.020226
$ws.Cells.Item(15, 15).Value = 0.06075245494223394
$ws.Cells.Item(15, 16).Value = 0.06075245494223393
$ws.Cells.Item(15, 17).Value = 0.003829588328666667
$ws.Cells.Item(15, 18).Value = 0.034466294958
$ws.Cells.Item(15, 19).Value = 0.0007874717360878301
$ws.Cells.Item(15, 20).Value = 0.0007874717360878301
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.011261
$ws.Cells.Item(16, 8).Value = 0.033783
$ws.Cells.Item(16, 9).Value = 0.01296197391260307
$ws.Cells.Item(16, 10).Value = 0.01296197391260308
$ws.Cells.Item(16, 13).Value = 0.5065093333333334
$ws.Cells.Item(16, 14).Value = 1.519528
$ws.Cells.Item(16, 15).Value = 0.0904849085922755
$ws.Cells.Item(16, 16).Value = 0.09048490859227548
$ws.Cells.Item(16, 17).Value = 0.005703801602666667
$ws.Cells.Item(16, 18).Value = 0.051334214424
$ws.Cells.Item(16, 19).Value = 0.001172863024657349
$ws.Cells.Item(16, 20).Value = 0.001172863024657349
$ws.Cells.Item(17, 7).Value = 0.7979533333333334
$ws.Cells.Item(17, 8).Value = 2.39386
$ws.Cells.Item(17, 9).Value = 0.9184841745974011
$ws.Cells.Item(17, 10).Value = 0.9184841745974011
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.569028
$ws.Cells.Item(17, 14).Value = 1.707084
$ws.Cells.Item(17, 15).Value = 0.1016535000995941
$ws.Cells.Item(17, 16).Value = 0.1016535000995941
$ws.Cells.Item(17, 17).Value = 0.45405778936
$ws.Cells.Item(17, 18).Value = 4.08652010424
$ws.Cells.Item(17, 19).Value = 0.0933671311339125
$ws.Cells.Item(17, 20).Value = 0.0933671311339125
$ws.Cells.Item(18, 7).Value = 0.7979533333333334
$ws.Cells.Item(18, 8).Value = 2.39386
$ws.Cells.Item(18, 9).Value = 0.9184841745974011
$ws.Cells.Item(18, 10).Value = 0.9184841745974011
$ws.Cells.Item(18, 15).Value = 0.1962512724671019
$ws.Cells.Item(18, 16).Value = 0.1962512724671019
$ws.Cells.Item(18, 17).Value = 0.8765996138666667
$ws.Cells.Item(18, 18).Value = 7.8893965248
$ws.Cells.Item(18, 19).Value = 0.1802536880056358
$ws.Cells.Item(18, 20).Value = 0.1802536880056358
$ws.Cells.Item(19, 7).Value = 0.7979533333333334
$ws.Cells.Item(19, 8).Value = 2.39386
$ws.Cells.Item(19, 9).Value = 0.9184841745974011
$ws.Cells.Item(19, 10).Value = 0.9184841745974011
$ws.Cells.Item(19, 13).Value = 3.083549
$ws.Cells.Item(19, 14).Value = 9.250647000000001
$ws.Cells.Item(19, 15).Value = 0.5508578638987945
$ws.Cells.Item(19, 16).Value = 0.5508578638987945
$ws.Cells.Item(19, 17).Value = 2.460528203046667
$ws.Cells.Item(19, 18).Value = 22.14475382742
$ws.Cells.Item(19, 19).Value = 0.5059542304435718
$ws.Cells.Item(19, 20).Value = 0.5059542304435718
$ws.Cells.Item(20, 7).Value = 0.7979533333333334
$ws.Cells.Item(20, 8).Value = 2.39386
$ws.Cells.Item(20, 9).Value = 0.9184841745974011
$ws.Cells.Item(20, 10).Value = 0.9184841745974011
$ws.Cells.Item(20, 13).Value = 0.3400753333333333
$ws.Cells.Item(20, 14).Value = 1.020226
$ws.Cells.Item(20, 15).Value = 0.06075245494223394
$ws.Cells.Item(20, 16).Value = 0.06075245494223393
$ws.Cells.Item(20, 17).Value = 0.2713642458177778
$ws.Cells.Item(20, 18).Value = 2.44227821236
$ws.Cells.Item(20, 19).Value = 0.05580016843238354
$ws.Cells.Item(20, 20).Value = 0.05580016843238354
$ws.Cells.Item(21, 7).Value = 0.7979533333333334
$ws.Cells.Item(21, 8).Value = 2.39386
$ws.Cells.Item(21, 9).Value = 0.9184841745974011
$ws.Cells.Item(21, 10).Value = 0.9184841745974011
$ws.Cells.Item(21, 13).Value = 0.5065093333333334
$ws.Cells.Item(21, 14).Value = 1.519528
$ws.Cells.Item(21, 15).Value = 0.0904849085922755
$ws.Cells.Item(21, 16).Value = 0.09048490859227548
$ws.Cells.Item(21, 17).Value = 0.4041708108977778
$ws.Cells.Item(21, 18).Value = 3.63753729808
$ws.Cells.Item(21, 19).Value = 0.08310895658189744
$ws.Cells.Item(21, 20).Value = 0.08310895658189743
